$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first review entry (row 2) - it was replaced by a later review
# that is now represented by the former row 3's content shifting up.
$ws.Rows("2:2").Delete()

# Update selection to B2 to mirror the saved selection state in the file.
$ws.Range("B2").Select()
